$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELSADCP Job Checklist")

# Update the date in A2
$ws.Range("A2").Value = "Date - 01.03.2024"

# Update the export dates in C14 (first 8 lines change from 31.01.24 to 28.02.24)
$c14 = "DGQ.R11KT6.BSIVM.TXT : last exported on 28.02.24`nDGQ.R11KT6.BSIVM.C.TXT : last exported on 28.02.24`nDGQ.R11KT6.BSIVM.US2.TXT : last exported on 28.02.24`nDGQ.R11KT6.PKAT.TXT : last exported on 28.02.24`nDGQ.R11KT6.PKAT.C.TXT : last exported on 28.02.24`nDGQ.R11KT6.PKAT.US2.TXT : last exported on 28.02.24`nDGQ.R11KT6.BSITK.TXT : last exported on 28.02.24`nDGQ.R11KT6.BSITX.TXT : last exported on 28.02.24`nDGQ.R31KT6.BSIVMC.TXT : last exported on 06.02.24`nDGQ.R31KT6.PKATC.TXT : last exported on 06.02.24`nDEL.KMPKT6.APOS.DATA.ZIP : last exported on 21.02.24"

$ws.Range("C14").Value = $c14
